$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 2, shifting the existing rows 2-4 down to 3-5.
$ws.Rows.Item(2).Insert()

# Populate the newly inserted row 2 with the "slug" identifiers that relate
# back to the human readable header row (row 1), so that the two rows can be
# linked together hierarchically (see issue #13).
$ws.Range("A2").Value = "municipio-superficie-medida"
$ws.Range("B2").Value = "municipio-zona-desfavorecida"
$ws.Range("C2").Value = "municipio-montana"
$ws.Range("D2").Value = "municipio-codigo"
$ws.Range("E2").Value = "municipio-nombre"
